$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 991060
$ws.Range("C4").Value = 1246.6
$ws.Range("D4").Value = 2442.3
$ws.Range("G4").Value = 359
$ws.Range("H4").Value = 1627
$ws.Range("I4").Value = 24791
$ws.Range("K4").Value = 1986.1
$ws.Range("L4").Value = 499
$ws.Range("M4").Value = 62.8
$ws.Range("N4").Value = 0.2
$ws.Range("O4").Value = 367143
$ws.Range("P4").Value = 461.8
$ws.Range("Q4").Value = 1562.7
$ws.Range("T4").Value = 43
$ws.Range("U4").Value = 300
$ws.Range("V4").Value = 20340
$ws.Range("X4").Value = 725.6
$ws.Range("Y4").Value = 506
$ws.Range("Z4").Value = 63.6
$ws.Range("AB4").Value = 494653
$ws.Range("AC4").Value = 622.2
$ws.Range("AD4").Value = 834.5
$ws.Range("AG4").Value = 323
$ws.Range("AH4").Value = 1016.5
$ws.Range("AI4").Value = 8295
$ws.Range("AK4").Value = 940.4
$ws.Range("AL4").Value = 526
$ws.Range("AM4").Value = 66.2
$ws.Range("AN4").Value = 0.7
$ws.Range("B5").Value = 1617872
$ws.Range("C5").Value = 675
$ws.Range("D5").Value = 2323.5
$ws.Range("G5").Value = 49
$ws.Range("H5").Value = 657
$ws.Range("I5").Value = 55081
$ws.Range("K5").Value = 1160.6
$ws.Range("L5").Value = 1394
$ws.Range("M5").Value = 58.2
$ws.Range("O5").Value = 579173
$ws.Range("P5").Value = 241.6
$ws.Range("Q5").Value = 1145.4
$ws.Range("T5").Value = 9
$ws.Range("U5").Value = 119
$ws.Range("V5").Value = 24618
$ws.Range("X5").Value = 406.2
$ws.Range("Y5").Value = 1426
$ws.Range("Z5").Value = 59.5
$ws.Range("AB5").Value = 902237
$ws.Range("AC5").Value = 376.4
$ws.Range("AD5").Value = 573.7
$ws.Range("AG5").Value = 126
$ws.Range("AH5").Value = 568
$ws.Range("AI5").Value = 5325
$ws.Range("AK5").Value = 590.9
$ws.Range("AL5").Value = 1527
$ws.Range("AM5").Value = 63.7
$ws.Range("AN5").Value = 0.3
$ws.Range("B6").Value = 188460
$ws.Range("C6").Value = 2048.5
$ws.Range("D6").Value = 2592.3
$ws.Range("F6").Value = 13
$ws.Range("G6").Value = 1172
$ws.Range("H6").Value = 3034.5
$ws.Range("I6").Value = 14651
$ws.Range("K6").Value = 2692.3
$ws.Range("L6").Value = 70
$ws.Range("M6").Value = 76.09999999999999
$ws.Range("O6").Value = 104009
$ws.Range("P6").Value = 1130.5
$ws.Range("Q6").Value = 5281.1
$ws.Range("S6").Value = 6
$ws.Range("T6").Value = 189
$ws.Range("U6").Value = 565.8
$ws.Range("V6").Value = 48717
$ws.Range("X6").Value = 1485.8
$ws.Range("Y6").Value = 70
$ws.Range("Z6").Value = 76.09999999999999
$ws.Range("AA6").Value = 1.1
$ws.Range("AB6").Value = 52618
$ws.Range("AC6").Value = 571.9
$ws.Range("AD6").Value = 647
$ws.Range("AG6").Value = 253.5
$ws.Range("AH6").Value = 1085.5
$ws.Range("AI6").Value = 2204
$ws.Range("AK6").Value = 877
$ws.Range("AL6").Value = 60
$ws.Range("AM6").Value = 65.2
$ws.Range("AN6").Value = 0.5
$ws.Range("B7").Value = 25556
$ws.Range("C7").Value = 464.7
$ws.Range("D7").Value = 1305.2
$ws.Range("H7").Value = 130
$ws.Range("I7").Value = 6144
$ws.Range("K7").Value = 1161.6
$ws.Range("O7").Value = 3929
$ws.Range("P7").Value = 71.40000000000001
$ws.Range("Q7").Value = 198.9
$ws.Range("U7").Value = 8.5
$ws.Range("V7").Value = 1010
$ws.Range("X7").Value = 187.1
$ws.Range("Y7").Value = 21
$ws.Range("Z7").Value = 38.2
$ws.Range("AA7").Value = -1.3
$ws.Range("AB7").Value = 15178
$ws.Range("AC7").Value = 276
$ws.Range("AD7").Value = 583.9
$ws.Range("AG7").Value = 6
$ws.Range("AH7").Value = 240
$ws.Range("AI7").Value = 2668
$ws.Range("AK7").Value = 523.4
$ws.Range("AL7").Value = 29
$ws.Range("AM7").Value = 52.7
$ws.Range("AN7").Value = -1.5
